$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Add the MS2 code path entry for esmith10laptop (column H = esmith10laptop,
# row 7 = MS2CodePath)
$ws.Range("H7").Value = "C:\E\GitHub\mRNADynamics"

# Keep the existing frozen/split pane boundary (between columns E and F) but
# scroll the right pane so it starts at column F, and leave the selection on
# the cell just entered (one column/row further, as in the saved workbook).
$aw = $excel.ActiveWindow
$aw.SplitColumn = 5
$aw.SplitRow = 0
$ws.Range("I7").Select()
